$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 963.7143
$ws.Range("K33").Value = 349.4
$ws.Range("I33").Value = 349.4
$ws.Range("M33").Value = -120.4
$ws.Range("L38").Value = 4350
$ws.Range("M38").Value = 230.000004
$ws.Range("J38").Value = 1450
$ws.Range("H38").Value = 247.71428
$ws.Range("K38").Value = 141.999996
$ws.Range("I38").Value = 47.333332
$ws.Range("N38").Value = -5094
$ws.Range("H41").Value = 322.77777
$ws.Range("K41").Value = 272.14285
$ws.Range("I41").Value = 272.14285
$ws.Range("M41").Value = 167.85715
$ws.Range("I92").Value = 2335.1428
$ws.Range("M92").Value = -1087.1428
$ws.Range("H92").Value = 2335.1428
$ws.Range("K92").Value = 2335.1428

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 837.1818
$ws.Range("M2").Value = -724.1818
$ws.Range("H2").Value = 1525.4286
$ws.Range("K2").Value = 837.1818
$ws.Range("L63").Value = 4500
$ws.Range("J63").Value = 4500
$ws.Range("H63").Value = 4500
$ws.Range("N63").Value = -5872
$ws.Range("J66").Value = 4500
$ws.Range("H66").Value = 4500
$ws.Range("N66").Value = -29364
$ws.Range("L66").Value = 22500
$ws.Range("H116").Value = 1525.4286
$ws.Range("K116").Value = 837.1818
$ws.Range("I116").Value = 837.1818
$ws.Range("M116").Value = 1456.8182

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 837.1818
$ws.Range("M3").Value = -723.1818
$ws.Range("K3").Value = 837.1818
$ws.Range("H3").Value = 1525.4286
$ws.Range("J20").Value = 2504.5
$ws.Range("H20").Value = 1739.75
$ws.Range("K20").Value = 975
$ws.Range("N20").Value = -2998.5
$ws.Range("L20").Value = 2504.5
$ws.Range("I20").Value = 975
$ws.Range("M20").Value = -728
$ws.Range("H82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H94").Value = 3527.6875
$ws.Range("K94").Value = 2716.4443
$ws.Range("N94").Value = -5472.7144
$ws.Range("L94").Value = 4570.7144
$ws.Range("I94").Value = 2716.4443
$ws.Range("M94").Value = -2265.4443
$ws.Range("J94").Value = 4570.7144
$ws.Range("H105").Value = 3765.2856
$ws.Range("K105").Value = 3142.8333
$ws.Range("N105").Value = -10994
$ws.Range("I105").Value = 3142.8333
$ws.Range("M105").Value = -1395.8333
$ws.Range("L105").Value = 7500
$ws.Range("J105").Value = 7500
$ws.Range("I107").Value = 748.4
$ws.Range("M107").Value = 1171.6
$ws.Range("L107").Value = 1900
$ws.Range("J107").Value = 1900
$ws.Range("N107").Value = -5740
$ws.Range("H107").Value = 940.3333
$ws.Range("K107").Value = 748.4

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 9000
$ws.Range("K41").Value = 9000
$ws.Range("I41").Value = 9000
$ws.Range("M41").Value = -8572
$ws.Range("J50").Value = 27141.143
$ws.Range("H50").Value = 27141.143
$ws.Range("N50").Value = -28391.143
$ws.Range("L50").Value = 27141.143
$ws.Range("N51").ClearContents()
$ws.Range("H51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("H59").Value = 28333.334
$ws.Range("N61").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("L71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("N88").Value = -21312
$ws.Range("L88").Value = 20500
$ws.Range("J88").Value = 20500
$ws.Range("H88").Value = 20500
$ws.Range("N91").Value = -23308
$ws.Range("H91").Value = 20500
$ws.Range("L91").Value = 20500
$ws.Range("J91").Value = 20500
$ws.Range("M92").Value = -20292
$ws.Range("L92").Value = 15300
$ws.Range("J92").Value = 15300
$ws.Range("H92").Value = 15300
$ws.Range("H96").Value = 12174.667
$ws.Range("N96").Value = -17666.667
$ws.Range("L96").Value = 12174.667
$ws.Range("J96").Value = 12174.667
$ws.Range("J106").Value = 80270.75
$ws.Range("H106").Value = 80270.75
$ws.Range("N106").Value = -82794.75
$ws.Range("L106").Value = 80270.75
$ws.Range("L107").Value = 621.75
$ws.Range("J107").Value = 621.75
$ws.Range("N107").Value = -4461.75
$ws.Range("H107").Value = 395.7
$ws.Range("H122").Value = 1261.7142
$ws.Range("K122").Value = 3785.1426
$ws.Range("I122").Value = 1261.7142
$ws.Range("M122").Value = -1335.1426
$ws.Range("I134").Value = 1503.8
$ws.Range("M134").Value = -1976.4
$ws.Range("H134").Value = 1482.5
$ws.Range("K134").Value = 4511.4
$ws.Range("L141").Value = 180266.67
$ws.Range("J141").Value = 180266.67
$ws.Range("M141").Value = -53820
$ws.Range("N141").Value = -190626.67
$ws.Range("H141").Value = 149950
$ws.Range("K141").Value = 59000
$ws.Range("I141").Value = 59000

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M7").Value = -443
$ws.Range("J7").Value = 296.33334
$ws.Range("N7").Value = -1113.00002
$ws.Range("H7").Value = 232.71428
$ws.Range("K7").Value = 555
$ws.Range("I7").Value = 185
$ws.Range("L7").Value = 889.0000200000001
$ws.Range("M36").Value = -1323.5
$ws.Range("H36").Value = 497.5
$ws.Range("K36").Value = 1492.5
$ws.Range("I36").Value = 497.5
$ws.Range("H92").Value = 793.75
$ws.Range("H117").Value = 3274.25
$ws.Range("N117").Value = -17480.9999
$ws.Range("L117").Value = 10596.9999
$ws.Range("J117").Value = 3532.3333

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J43").Value = 8333.333000000001
$ws.Range("N43").Value = -8635.333000000001
$ws.Range("H43").Value = 9648.058999999999
$ws.Range("K43").Value = 11127.125
$ws.Range("I43").Value = 11127.125
$ws.Range("M43").Value = -10976.125
$ws.Range("L43").Value = 8333.333000000001
$ws.Range("K70").Value = 9750.5
$ws.Range("N70").ClearContents()
$ws.Range("H70").Value = 9750.5
$ws.Range("I70").Value = 9750.5
$ws.Range("M70").Value = -9480.5
$ws.Range("L70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K73").Value = 9750.5
$ws.Range("N73").ClearContents()
$ws.Range("H73").Value = 9750.5
$ws.Range("I73").Value = 9750.5
$ws.Range("M73").Value = -8814.5
$ws.Range("L73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("H80").Value = 8253
$ws.Range("K80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 8253
$ws.Range("K83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("J122").Value = 1599.5
$ws.Range("H122").Value = 1644.3334
$ws.Range("K122").Value = 4971.428400000001
$ws.Range("N122").Value = -9698.5
$ws.Range("L122").Value = 4798.5
$ws.Range("I122").Value = 1657.1428
$ws.Range("M122").Value = -2521.428400000001
$ws.Range("I132").Value = 3743.8125
$ws.Range("M132").Value = -8701.4375
$ws.Range("H132").Value = 3697.7646
$ws.Range("K132").Value = 11231.4375

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -1588
$ws.Range("H7").Value = 1700
$ws.Range("K7").Value = 1700
$ws.Range("I7").Value = 1700
$ws.Range("H126").Value = 1700
$ws.Range("K126").Value = 5100
$ws.Range("I126").Value = 1700
$ws.Range("M126").Value = -2630
$ws.Range("I132").Value = 3833.3333
$ws.Range("M132").Value = -8969.999899999999
$ws.Range("H132").Value = 3833.3333
$ws.Range("K132").Value = 11499.9999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 7746655.5
$ws.Range("K100").Value = 23235382
$ws.Range("N100").Value = -10248.667
$ws.Range("L100").Value = 9166.666999999999
$ws.Range("I100").Value = 11617691
$ws.Range("M100").Value = -23234841
$ws.Range("J100").Value = 4583.3335
$ws.Range("I132").Value = 743.2143
$ws.Range("M132").Value = 300.3571000000002
$ws.Range("H132").Value = 760.3333
$ws.Range("K132").Value = 2229.6429
